$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E retain their original Text format so that
# numeric-looking values (e.g. "0.999") are written as text, not numbers,
# matching the inlineStr cell type used throughout the sheet.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.124.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.546.68"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.88%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "614.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.60"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.37%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.25%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.10%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.91%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.55"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.110.82"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "632.13"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +10.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.97"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.178.12"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.66%  "

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.571.17"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.95"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.98%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.996"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.60"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.73"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.96"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +8.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.90"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.82%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.02"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.62"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +7.17%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.46%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.06"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.36"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.21%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "63.88"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.68"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +17.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.24"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "530.03"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.401"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.38%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.22"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0782"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.52"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.508.15"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.49%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.09%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.99%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.145"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.30%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.15"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.38%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.41"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.26"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.30%  "
